$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying data for row 2 and row 3 got swapped for columns
# D (Fecha), M (Volumen), N (Precio mínimo), O (Precio máximo),
# P (Precio promedio ponderado), R (Origen) and S (Precio $/Kg).

# Row 2 new values (previously held by row 3)
$ws.Range("D2").Value = 44235
$ws.Range("M2").Value = 70
$ws.Range("N2").Value = 42000
$ws.Range("O2").Value = 42000
$ws.Range("P2").Value = 42000
$ws.Range("R2").Value = "Región de Arica y Parinacota"
$ws.Range("S2").Value = 2333

# Row 3 new values (previously held by row 2)
$ws.Range("D3").Value = 44417
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 26000
$ws.Range("O3").Value = 26000
$ws.Range("P3").Value = 26000
$ws.Range("R3").Value = "Perú"
$ws.Range("S3").Value = 1444
